$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('G2').Value2 = 'system, backup@backdoor.com, System'
$ws.Range('G4').Value2 = 'backup@backdoor.com, System'
$ws.Range('G5').Value2 = 'backup@backdoor.com, System'
$ws.Range('G7').Value2 = 'admin@admin.com, System'
$ws.Range('G8').Value2 = 'backup@backdoor.com, System'
$ws.Range('G11').Value2 = 'dnasr281@gmail.com, System'
$ws.Range('G17').Value2 = 'dnasr281@gmail.com, System'
$ws.Range('G28').Value2 = 'system, backup@backdoor.com, System'
$ws.Range('G30').Value2 = 'backup@backdoor.com, System'
$ws.Range('G31').Value2 = 'backup@backdoor.com, System'
$ws.Range('G33').Value2 = 'admin@admin.com, System'
$ws.Range('G34').Value2 = 'backup@backdoor.com, System'
$ws.Range('G37').Value2 = 'dnasr281@gmail.com, System'
$ws.Range('G43').Value2 = 'dnasr281@gmail.com, System'
$ws.Range('G54').Value2 = 'system, backup@backdoor.com, System'
$ws.Range('G56').Value2 = 'backup@backdoor.com, System'
$ws.Range('G57').Value2 = 'backup@backdoor.com, System'
$ws.Range('G59').Value2 = 'admin@admin.com, System'
$ws.Range('G60').Value2 = 'backup@backdoor.com, System'
$ws.Range('G63').Value2 = 'dnasr281@gmail.com, System'
$ws.Range('G69').Value2 = 'dnasr281@gmail.com, System'
$ws.Range('G80').Value2 = 'backup@backdoor.com, System'
$ws.Range('G81').Value2 = 'backup@backdoor.com, System'
$ws.Range('G82').Value2 = 'backup@backdoor.com, System'
$ws.Range('G93').Value2 = 'dnasr281@gmail.com, System'
$ws.Range('G94').Value2 = 'dnasr281@gmail.com, System'
$ws.Range('G96').Value2 = 'dnasr281@gmail.com, System'
$ws.Range('G106').Value2 = 'backup@backdoor.com, System'
$ws.Range('G107').Value2 = 'backup@backdoor.com, System'
$ws.Range('G108').Value2 = 'backup@backdoor.com, System'
$ws.Range('G119').Value2 = 'dnasr281@gmail.com, System'
$ws.Range('G120').Value2 = 'dnasr281@gmail.com, System'
$ws.Range('G122').Value2 = 'dnasr281@gmail.com, System'
$ws.Range('G132').Value2 = 'backup@backdoor.com, System'
$ws.Range('G133').Value2 = 'backup@backdoor.com, System'
$ws.Range('G134').Value2 = 'backup@backdoor.com, System'
$ws.Range('G145').Value2 = 'dnasr281@gmail.com, System'
$ws.Range('G146').Value2 = 'dnasr281@gmail.com, System'
$ws.Range('G148').Value2 = 'dnasr281@gmail.com, System'
